# Getana Deliverable 3 Product Backlog - "added story 33 to backlog"
#
# A new user story (ID 33 - review privacy policy/implications statement)
# is inserted as a new row 23 in the "Sprint 4" backlog sheet, pushing the
# existing L-priority backlog rows (previously rows 23-33) down to rows
# 24-34 and renumbering their "L, NN" priority labels by +1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 23 - it inherits formatting from the row above (row 22),
# matching Excel's default "shift cells down" row-insert behavior.
$ws.Rows("23:23").Insert()

# Populate the new backlog item.
$ws.Range("A23").Value = 33
$ws.Range("B23").Value = "As a user, I want to be able to review the privacy policy/implications statement at any time, preferably from within the settings menu, so that do not have to take special measures to review the statement if I have forgotten its contents."
$ws.Range("C23").Value = "undecided"
$ws.Range("D23").Value = "L, 13"
$ws.Range("E23").Value = "W"
$ws.Range("F23").Value = 1

# The rows that were pushed down keep their own text, but their priority
# labels ("L, 13" .. "L, 23") need to be bumped by one ("L, 14" .. "L, 24")
# since the new story claimed "L, 13".
$ws.Range("D24").Value = "L, 14"
$ws.Range("D25").Value = "L, 15"
$ws.Range("D26").Value = "L, 16"
$ws.Range("D27").Value = "L, 17"
$ws.Range("D28").Value = "L, 18"
$ws.Range("D29").Value = "L, 19"
$ws.Range("D30").Value = "L, 20"
$ws.Range("D31").Value = "L, 21"
$ws.Range("D32").Value = "L, 22"
$ws.Range("D33").Value = "L, 23"
$ws.Range("D34").Value = "L, 24"

# Match the author's final cursor position.
$ws.Range("D36").Select()
